# [ADD] Column and show 'cannot be retrieved' status
#
# Adds 4 new header columns to the "Budget Asset Report" sheet's header row
# (row 9): "Division", "Section Code", "Section Name" (inserted right after
# "Org"), and "Status" (inserted right after "Investment Asset Name").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert 3 new columns right after "Org" (old column B), -------
# pushing "Investment Asset Code" / "Investment Asset Name" and everything
# after them three columns to the right (old B,C -> E,F).
$ws.Range("B:D").EntireColumn.Insert()

# --- Step 2: insert 1 new column right after "Investment Asset Name" ------
# (now column F), pushing "Budget Plan" and everything after it one column
# to the right (old D.. -> H..).
$ws.Range("G:G").EntireColumn.Insert()

# Inserting whole columns also stamps the (blank) rows above the header
# (rows 1-7, which only ever had data in column A) with blank, formatted
# cells in the new columns. Those rows only ever had an "A" cell, so clear
# the spurious blanks back out again.
$ws.Range("B1:D7").Clear()
$ws.Range("G1:G7").Clear()

# --- Step 3: fill in the new header cells on row 9 ------------------------
$ws.Range("B9").Value2 = "Division"
$ws.Range("C9").Value2 = "Section Code"
$ws.Range("D9").Value2 = "Section Name"
$ws.Range("G9").Value2 = "Status"

# Widen "Division" (B) to fit its header text, and give the other brand new
# columns (C, D "Section Code"/"Section Name" and G "Status") the same
# width as their neighbouring column groups.
$ws.Range("B1").EntireColumn.ColumnWidth = 33.07
$ws.Range("C1:D1").EntireColumn.ColumnWidth = 16.39
$ws.Range("G1").EntireColumn.ColumnWidth = 14.72

# Copy the existing header style (row 9 headers use the same style, e.g. A9)
# onto the freshly-inserted header cells so they render identically to the
# rest of the header row.
$ws.Range("A9").Copy()
$headerCells = $ws.Range("B9,C9,D9,G9")
$headerCells.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- Step 4: drop the stray legacy cells that used to trail off row 9 at ---
# the far right of the sheet (columns AMG:AMJ originally, now shifted to
# AMK:AMN after the inserts above) so the sheet's used range shrinks back
# down to A1:S9 instead of stretching out to the last column.
$ws.Range("AMK9:AMN9").ClearContents()
$ws.Range("AMK9:AMN9").ClearFormats()

# --- Step 5: update the active selection / view ----------------------------
$ws.Range("O9").Select()
